$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2 through 15
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22)
$ws.Range("C2:C15").Value = 45221
